$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = -0.388484531405109
$ws.Range("B2").Value = -0.3914174089723574
$ws.Range("A3").Value = -0.4244740690991697
$ws.Range("B3").Value = -0.3669978345692769
$ws.Range("A4").Value = -0.3835623551530521
$ws.Range("B4").Value = -0.3953804278845049
$ws.Range("A5").Value = -0.2128215971196397
$ws.Range("B5").Value = -0.2768032666367231
$ws.Range("A6").Value = -0.1178138514474628
$ws.Range("B6").Value = -0.09251671575425606
$ws.Range("A7").Value = -0.3675778940467194
$ws.Range("B7").Value = -0.3540105305702276
$ws.Range("A8").Value = -0.4959516217926717
$ws.Range("B8").Value = -0.4014203587121239
$ws.Range("A9").Value = -0.5074432674552886
$ws.Range("B9").Value = -0.4198084368880229
$ws.Range("A10").Value = -0.5631564098676105
$ws.Range("B10").Value = -0.3666393816881962
$ws.Range("A11").Value = -0.2264848141108131
$ws.Range("B11").Value = -0.3058016753278167
$ws.Range("A12").Value = -0.3131388064915851
$ws.Range("B12").Value = -0.2953366344957017
$ws.Range("A13").Value = 0.1123417641906094
$ws.Range("B13").Value = 0.007297919376174351
$ws.Range("A14").Value = -0.1039719205296165
$ws.Range("B14").Value = -0.1276674254616423
$ws.Range("A15").Value = -0.1203792822528296
$ws.Range("B15").Value = -0.04764335489635556
$ws.Range("A16").Value = -0.1991460628533611
$ws.Range("B16").Value = -0.1067347119136941
$ws.Range("A17").Value = 0.01623643607608446
$ws.Range("B17").Value = 0.08112199228031655
$ws.Range("A18").Value = 0.08781708433588127
$ws.Range("B18").Value = 0.09922421118552731
$ws.Range("A19").Value = 0.004903479896118251
$ws.Range("B19").Value = 0.06325074652347444
$ws.Range("A20").Value = -0.06120257041936774
$ws.Range("B20").Value = -0.02152986183189824
$ws.Range("A21").Value = -0.1233175890892114
$ws.Range("B21").Value = -0.07416718027996146
$ws.Range("A22").Value = 0.03556270798435462
$ws.Range("B22").Value = 0.09065920898830772
$ws.Range("A23").Value = 0.130141022801464
$ws.Range("B23").Value = 0.07021733751151711
$ws.Range("A24").Value = 0.5287553875309271
$ws.Range("B24").Value = 0.3487956199618114
$ws.Range("A25").Value = 0.2038741664190011
$ws.Range("B25").Value = 0.1534144205454203
$ws.Range("A26").Value = 0.18729772488959
$ws.Range("B26").Value = 0.1555936490986327
$ws.Range("A27").Value = 0.1498935711103744
$ws.Range("B27").Value = 0.1404457159736983
$ws.Range("A28").Value = 0.2773270094443454
$ws.Range("B28").Value = 0.1828995105258386
$ws.Range("A29").Value = 0.5388910676282144
$ws.Range("B29").Value = 0.4122407233067132
$ws.Range("A30").Value = 0.2225758275239678
$ws.Range("B30").Value = 0.1825663908112712
$ws.Range("A31").Value = 0.1350011835774995
$ws.Range("B31").Value = 0.110429181616766
$ws.Range("A32").Value = 0.1658953502344915
$ws.Range("B32").Value = 0.1633259084726907
$ws.Range("A33").Value = 0.1513595288739798
$ws.Range("B33").Value = 0.1521463839013246
$ws.Range("A34").Value = 0.1653742679550473
$ws.Range("B34").Value = 0.09916269753708935
$ws.Range("A35").Value = 0.177073582281517
$ws.Range("B35").Value = 0.1143981581327905
$ws.Range("A36").Value = 0.1127578816511595
$ws.Range("B36").Value = 0.06205583585826165
$ws.Range("A37").Value = 0.1980261628450812
$ws.Range("B37").Value = 0.05690202298089145
$ws.Range("A38").Value = 0.3637670641608936
$ws.Range("B38").Value = 0.2984508649202887
$ws.Range("A39").Value = -0.04301544970390719
$ws.Range("B39").Value = -0.1390648837184182
$ws.Range("A40").Value = 0.1965181960545027
$ws.Range("B40").Value = 0.1418565862980938
$ws.Range("A41").Value = -0.01445158207321362
$ws.Range("B41").Value = -0.03275906288244444
$ws.Range("A42").Value = 0.1553468126097762
$ws.Range("B42").Value = 0.1522238681624847
$ws.Range("A43").Value = 0.2597195569870762
$ws.Range("B43").Value = 0.2004962602694975
$ws.Range("A44").Value = -0.06259852868307275
$ws.Range("B44").Value = -0.08756464595435048
$ws.Range("A45").Value = -0.1398700938403052
$ws.Range("B45").Value = -0.1232395697573842
$ws.Range("A46").Value = -0.1768612919663619
$ws.Range("B46").Value = -0.1807300489219615
$ws.Range("A47").Value = -0.178646714833852
$ws.Range("B47").Value = -0.1835581835252835
$ws.Range("A48").Value = -0.2200293231343671
$ws.Range("B48").Value = -0.2049292595153646
$ws.Range("A49").Value = -0.2173539483118564
$ws.Range("B49").Value = -0.2094446973979209
$ws.Range("A50").Value = -0.152842649386205
$ws.Range("B50").Value = -0.1514611093938477
$ws.Range("A51").Value = -0.2088565662868457
$ws.Range("B51").Value = -0.2291447312577234
$ws.Range("A52").Value = -0.2088565662868457
$ws.Range("B52").Value = -0.2291447312577234
$ws.Range("A53").Value = -0.1871987613209684
$ws.Range("B53").Value = -0.17867393589209
$ws.Range("A54").Value = -0.2292806647691895
$ws.Range("B54").Value = -0.2193613656124554
$ws.Range("A55").Value = -0.1736771729705233
$ws.Range("B55").Value = -0.1603403775775097
$ws.Range("A56").Value = -0.1693028641325004
$ws.Range("B56").Value = -0.1727691167437934
$ws.Range("A57").Value = -0.2114115535736886
$ws.Range("B57").Value = -0.1766141786474135
$ws.Range("A58").Value = -0.1970921436752836
$ws.Range("B58").Value = -0.2233805576375427
$ws.Range("A59").Value = -0.2571481345669195
$ws.Range("B59").Value = -0.2520434188710283
$ws.Range("A60").Value = -0.3057288860745265
$ws.Range("B60").Value = -0.285298308011424
$ws.Range("A61").Value = -0.2299735983039812
$ws.Range("B61").Value = -0.196952230324226
$ws.Range("A62").Value = -0.219350464163844
$ws.Range("B62").Value = -0.134454463472991
$ws.Range("A63").Value = -0.3594398267956446
$ws.Range("B63").Value = -0.3961199357731875
$ws.Range("A64").Value = -0.2714261672907615
$ws.Range("B64").Value = -0.2569065801217919
$ws.Range("A65").Value = -0.3511372033851052
$ws.Range("B65").Value = -0.321417880495291
$ws.Range("A66").Value = -0.1712564983919799
$ws.Range("B66").Value = -0.137149999642194
$ws.Range("A67").Value = -0.1685632364568675
$ws.Range("B67").Value = -0.1502495909086008
